# Applies the "Add files via upload" revision to the
# "Commitment of Traders" workbook:
#   - Header cell H3 ("Net") is relabelled "Speculator - Hedger"
#     (a brand-new shared string is created for this).
#   - The merged header cells H3:H4 pick up word-wrap so the longer
#     caption can wrap onto multiple lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header in H3 (merged with H4) from "Net" to "Speculator - Hedger".
$ws.Range("H3").Value() = "Speculator - Hedger"

# Turn on wrap-text for the merged header cell so the new, longer
# caption fits nicely (keeps the existing center / vertical-center alignment).
$ws.Range("H3:H4").WrapText = $true
